$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

$range = $ws.Range("C2:C$lastRow")
$range.Value = (Get-Date -Year 2026 -Month 2 -Day 22 -Hour 0 -Minute 0 -Second 0).Date
